$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts Connect Subscription ID and
# everything to its right one column to the right).
$ws.Columns("C").Insert()

# Populate the new column's header.
$ws.Range("C1").Value = "Assignee Name"

# New column should match the width of the neighboring A:B columns.
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Match the active cell shown in the diff.
$ws.Range("C1").Select()
